$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("I14").Value = 3
$ws.Range("N14").Value = -62.5
$ws.Range("C15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -50
$ws.Range("E16").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -10
$ws.Range("M15").Value = 800
$ws.Range("N15").Value = 12.5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = -61.538461538461
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 64
$ws.Range("K16").Value = -1.5625
$ws.Range("L16").Value = -11.267605633802
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = -85
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 133.333333333333
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 63
$ws.Range("K17").Value = 63.492063492063
$ws.Range("L17").Value = 66.129032258064
$ws.Range("M17").Value = 255.172413793103
$ws.Range("N17").Value = 30.379746835443
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 64.285714285714
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = 24.742268041237
$ws.Range("L18").Value = -10.37037037037
$ws.Range("M18").Value = 39.080459770114
$ws.Range("N18").Value = -71.462264150943
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 35
$ws.Range("F19").Value = 103
$ws.Range("G19").Value = 101
$ws.Range("H19").Value = 1.980198019801
$ws.Range("I19").Value = 551
$ws.Range("J19").Value = 575
$ws.Range("K19").Value = -4.173913043478
$ws.Range("L19").Value = -8.319467554076
$ws.Range("M19").Value = 5.353728489483
$ws.Range("N19").Value = -69.675288937809
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 9
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -38.235294117647
$ws.Range("L20").Value = -40
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -95.24886877828
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 18.75
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = 9.395973154362
$ws.Range("I21").Value = 871
$ws.Range("J21").Value = 842
$ws.Range("K21").Value = 3.444180522565
$ws.Range("L21").Value = -4.704595185995
$ws.Range("M21").Value = 26.231884057971
$ws.Range("N21").Value = -72.764227642276
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 43
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = 4.878048780487
$ws.Range("L22").Value = -21.818181818181
$ws.Range("M22").Value = 16.216216216216
$ws.Range("C24").Value = 78
$ws.Range("D24").Value = 74
$ws.Range("E24").Value = 5.405405405405
$ws.Range("F24").Value = 270
$ws.Range("G24").Value = 293
$ws.Range("H24").Value = -7.849829351535
$ws.Range("I24").Value = 1941
$ws.Range("J24").Value = 1887
$ws.Range("K24").Value = 2.861685214626
$ws.Range("L24").Value = -1.371951219512
$ws.Range("M24").Value = 126.223776223776
$ws.Range("C25").Value = 78
$ws.Range("E25").Value = 5.405405405405
$ws.Range("F25").Value = 252
$ws.Range("G25").Value = 314
$ws.Range("H25").Value = -19.745222929936
$ws.Range("I25").Value = 1907
$ws.Range("J25").Value = 1912
$ws.Range("K25").Value = -0.26150627615
$ws.Range("L25").Value = -4.697651174412
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 22.222222222222
$ws.Range("F26").Value = 39
$ws.Range("H26").Value = 39.285714285714
$ws.Range("I26").Value = 232
$ws.Range("J26").Value = 181
$ws.Range("K26").Value = 28.17679558011
$ws.Range("L26").Value = 39.759036144578
$ws.Range("M26").Value = 82.67716535433
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -50
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -9.090909090909
$ws.Range("L27").Value = -28.571428571428
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 37.5
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = 20.454545454545
$ws.Range("L28").Value = 1.923076923076
$ws.Range("C29").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("I29").Value = 3
$ws.Range("L29").Value = 200
$ws.Range("N29").Value = -25
$ws.Range("C30").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("I30").Value = 3
$ws.Range("L30").Value = 200
$ws.Range("N30").Value = -25
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# --- Header text updates (rich text substring edits) ---
$ws.Range("A8").Characters(21, 2).Text = "26"
$ws.Range("C9").Characters(27, 9).Text = "6/24/2024"
$ws.Range("C9").Characters(47, 9).Text = "6/30/2024"

# --- Column H width (bestFit narrowed after data changed to shorter strings) ---
$ws.Columns.Item(8).ColumnWidth = 5.43
